# Updated sampling ranges for experiment and added climate change factor to hydropower.
#
# The "General" subsector table on the first worksheet (strategy_id-0) gets a
# new row inserted right above "elasticity_gnrl_rate_occupancy_to_gdppc"
# (currently row 4), pushing it and every row below it down by one. The new
# row 4 holds the "climate_change_factor_gnrl_hydropower_availability"
# variable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 4; existing rows 4-11 shift to 5-12.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the climate_change_factor_gnrl_hydropower_availability entry.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5
$ws.Range("J4:AS4").Value = 1
